$d = $word.ActiveDocument

# wdHeaderFooterIndex constants
$wdHeaderFooterPrimary   = 1
$wdHeaderFooterFirstPage = 2
$wdHeaderFooterEvenPages = 3

foreach ($sec in $d.Sections) {

    # --- Headers: BTEC logo picture (header "first page" story) goes from
    #     image1.jpg -> image2.jpg ---
    foreach ($idx in @($wdHeaderFooterPrimary, $wdHeaderFooterFirstPage, $wdHeaderFooterEvenPages)) {
        $hdr = $sec.Headers.Item($idx)
        if ($hdr.Exists) {
            $shapes = $hdr.Range.InlineShapes
            for ($i = 1; $i -le $shapes.Count; $i++) {
                $shp = $shapes.Item($i)
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    $shp.Name = "image2.jpg"
                }
            }
        }
    }

    # --- Footers: Pearson Edexcel logo pictures go from image2.png -> image1.png ---
    foreach ($idx in @($wdHeaderFooterPrimary, $wdHeaderFooterFirstPage, $wdHeaderFooterEvenPages)) {
        $ftr = $sec.Footers.Item($idx)
        if ($ftr.Exists) {
            $shapes = $ftr.Range.InlineShapes
            for ($i = 1; $i -le $shapes.Count; $i++) {
                $shp = $shapes.Item($i)
                if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                    $shp.Name = "image1.png"
                }
            }
        }
    }
}
